$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Programa de hub Python"
$ws.Range("B15").Value = "Programa cerradura"
$ws.Range("B16").Value = "Hearbeats"
$ws.Range("B17").Value = "Servicios REST"
$ws.Range("B18").Value = "Pruebas postman"
